# Update Name of Algo
# Apply updated RandomForest imputation values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -11.5732
$ws.Range("B7").Value = 4.870099999999997
$ws.Range("A10").Value = -21.75489999999999
$ws.Range("E10").Value = 16.12309999999999
$ws.Range("A12").Value = -21.55350000000001
$ws.Range("E14").Value = 16.84560000000001
$ws.Range("B15").Value = 4.506299999999995
$ws.Range("A18").Value = -22.08260000000001
$ws.Range("C18").Value = -11.6029
$ws.Range("C19").Value = -11.5732
$ws.Range("B20").Value = 9.389199999999995
$ws.Range("C27").Value = -12.8437
$ws.Range("B29").Value = 5.084400000000005
$ws.Range("B30").Value = 4.495200000000001
$ws.Range("B31").Value = 4.816899999999999
$ws.Range("E32").Value = 16.58449999999999
$ws.Range("E35").Value = 16.33910000000001
$ws.Range("A37").Value = -19.53359999999999
$ws.Range("B40").Value = 9.258699999999994
$ws.Range("C42").Value = -11.70570000000001
$ws.Range("E43").Value = 17.13490000000001
$ws.Range("C44").Value = -13.86619999999998
$ws.Range("C47").Value = -12.0235
$ws.Range("E49").Value = 15.54420000000001
$ws.Range("A55").Value = -22.42900000000001
$ws.Range("E56").Value = 16.6721
$ws.Range("C58").Value = -11.84929999999999
$ws.Range("A68").Value = -21.49460000000001
$ws.Range("B68").Value = 4.6753
$ws.Range("E69").Value = 17.38130000000002
$ws.Range("C73").Value = -12.48190000000001
$ws.Range("B76").Value = 5.639400000000002
$ws.Range("A77").Value = -20.32459999999998
$ws.Range("A78").Value = -20.27409999999998
$ws.Range("E81").Value = 16.81109999999999
$ws.Range("B87").Value = 5.767299999999995
$ws.Range("B88").Value = 5.043799999999997
$ws.Range("E92").Value = 18.37700000000002
$ws.Range("C95").Value = -12.2697
$ws.Range("B96").Value = 5.326600000000006
$ws.Range("B98").Value = 6.744800000000001
$ws.Range("B101").Value = 9.226199999999993
$ws.Range("C101").Value = -12.91610000000001
$ws.Range("B102").Value = 8.624500000000008
